$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 179246
$ws.Range("C4").Value = 169202
$ws.Range("C7").Value = 5.6
$ws.Range("C8").Value = 65.05
